$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 currently holds the text "R40"; change it to the text "1".
#
# A leading apostrophe forces Excel to store a numeric-looking entry as
# text (t="s", a new shared-string entry) instead of silently coercing it
# to a number - but doing that also flips on the cell's "stored as text"
# quote-prefix flag, which would stamp B11 with a brand new style. Since
# B11's existing style must stay exactly as it was, stash its formatting
# in a scratch cell first and paste it straight back once the text value
# is in place.
$cell = $ws.Range("B11")
$scratch = $ws.Range("Z1")

$cell.Copy()
$scratch.PasteSpecial(-4122)  # xlPasteFormats

$cell.Value = "'1"

$scratch.Copy()
$cell.PasteSpecial(-4122)     # xlPasteFormats - restores B11's original style
$scratch.Clear()              # remove the scratch cell, leaving no trace
